$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A146").Value = "IMX-USD"
$ws.Range("A147").Value = "GRT-USD"
